# RNASeq_MiniLecture_03_03_DifferentialExpression.pptx
# "update names and dates"
#
# 1) Update the author/date credit box on slide 1 (shape "Title 1", the
#    second shape with that name) with the new list of contributors and
#    the new workshop date range.
# 2) Update the notes-master "updated automatically" date field's cached
#    text (best effort - PowerPoint date fields are normally refreshed
#    from the system clock, not hand-edited).

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 1: credits / date textbox
# -----------------------------------------------------------------
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(6)
$tr = $shape.TextFrame.TextRange

# ---- Paragraph 1: "Kelsy Cotto, Arpad Danos, Harriet Dashnow, Felicia
#      Gomez, Sharon Freshour, Obi Griffith, " -----------------------
$para1 = $tr.Characters(1, 78)
$para1.Text = "Kelsy Cotto, Arpad Danos, Harriet Dashnow, Felicia Gomez, Sharon Freshour, Obi Griffith, "

# Re-touch (no visual change - same font re-applied) the three surname
# runs so they remain separate runs from the surrounding text, matching
# how PowerPoint splits runs around words flagged by the spell checker.
$tr.Characters(20, 5).Font.Name = "Calibri"   # Danos
$tr.Characters(35, 7).Font.Name = "Calibri"   # Dashnow
$tr.Characters(66, 8).Font.Name = "Calibri"   # Freshour

# ---- Paragraph 2: "Malachi Griffith, Jason Kunisaki, Chris Miller,
#      Jonathan Preall, Aaron Quinlan" --------------------------------
$para2Start = 91
$para2 = $tr.Characters($para2Start, 71)
$para2.Text = "Malachi Griffith, Jason Kunisaki, Chris Miller, Jonathan Preall, Aaron Quinlan"

$tr.Characters($para2Start + 24, 8).Font.Name = "Calibri"   # Kunisaki
$tr.Characters($para2Start + 57, 6).Font.Name = "Calibri"   # Preall

# ---- Paragraph 3: refresh the workshop date ----
$para3Start = $para2Start + 71 + 8  # start of paragraph 3 after the resize of paragraph 2
$dateRun = $tr.Characters($para3Start + 59, 21)
$dateRun.Text = "November 11-19, 2021"

# -----------------------------------------------------------------
# Notes master: cached "updated automatically" date field
# -----------------------------------------------------------------
$nm = $p.NotesMaster
$dateShape = $nm.Shapes.Item(2)
$dateShape.TextFrame.TextRange.Text = "11/10/21"
